$d = $word.ActiveDocument

# 1. Remove the stray _GoBack bookmark that currently sits on the first
#    paragraph ("به نام خدا").
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. The title paragraph currently reads "تمرین اول شبکه های کامپیوتری"
#    ("Computer Networks homework one"). It becomes homework four, and a
#    new _GoBack bookmark is left right after the inserted word "چهارم".
#    Locate that specific paragraph by its own text (rather than a fixed
#    index) so we never accidentally touch the later "سوال اول"
#    ("Question one") paragraph, which also contains the word "اول".
$titlePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*تمرین*اول*شبکه*") {
        $titlePara = $p
        break
    }
}
if ($null -eq $titlePara) {
    throw "Could not find the homework-title paragraph"
}

$searchRange = $titlePara.Range
$found = $searchRange.Find.Execute("اول", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the word to replace in the title paragraph"
}

$splitPoint = $searchRange.Start
$searchRange.Text = "چهارم"
$afterWord = $searchRange.End

# Re-insert the _GoBack bookmark immediately after "چهارم" (collapsed
# range); inserting a bookmark there also splits the run at that point.
$bmRange = $d.Range($afterWord, $afterWord)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Force the run boundary between "تمرین " and "چهارم" the same way, via a
# throwaway bookmark that gets removed immediately afterwards.
$splitRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("zzTmpSplit", $splitRange)
$d.Bookmarks("zzTmpSplit").Delete()
